# Update the "want-to-go" counters (column F) for a handful of events across
# three sheets (展览 / 演出 / 全部类型), matching the regenerated scrape output.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 1053
$ws.Range("F4").Value  = 1461
$ws.Range("F8").Value  = 597
$ws.Range("F16").Value = 1048
$ws.Range("F20").Value = 4257
$ws.Range("F22").Value = 3289
$ws.Range("F25").Value = 3137
$ws.Range("F30").Value = 3100
$ws.Range("F36").Value = 1121
$ws.Range("F43").Value = 484
$ws.Range("F49").Value = 3691

# Sheet "演出" (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 982

# Sheet "全部类型" (All categories)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 1960
$ws.Range("F5").Value  = 1461
$ws.Range("F13").Value = 982
$ws.Range("F17").Value = 1048
$ws.Range("F21").Value = 4257
$ws.Range("F25").Value = 3289
$ws.Range("F26").Value = 3137
$ws.Range("F29").Value = 3100
$ws.Range("F41").Value = 484
$ws.Range("F49").Value = 3691
